$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 208
$ws.Range("F4").Value = 389
$ws.Range("F6").Value = 5
$ws.Range("F7").Value = 552
$ws.Range("F9").Value = 9779
$ws.Range("F10").Value = 51
$ws.Range("F11").Value = 2643
$ws.Range("F12").Value = 206
$ws.Range("F13").Value = 2388
$ws.Range("F14").Value = 2655
$ws.Range("F15").Value = 1397
$ws.Range("F17").Value = 2076
$ws.Range("F19").Value = 79
$ws.Range("F20").Value = 366
$ws.Range("F22").Value = 70
$ws.Range("F23").Value = 297
$ws.Range("F24").Value = 61
$ws.Range("F25").Value = 147
$ws.Range("F27").Value = 1286
$ws.Range("F29").Value = 95
$ws.Range("F30").Value = 122
$ws.Range("F31").Value = 243
$ws.Range("F32").Value = 1676
$ws.Range("F33").Value = 2807
$ws.Range("F34").Value = 27
$ws.Range("F35").Value = 989
$ws.Range("F36").Value = 358
$ws.Range("F37").Value = 1
$ws.Range("F38").Value = 1275
$ws.Range("F39").Value = 49
$ws.Range("F40").Value = 54
$ws.Range("F41").Value = 50

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F13").Value = 34
$ws.Range("F14").Value = 155

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 951
$ws.Range("F5").Value = 1635

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 208
$ws.Range("F4").Value = 951
$ws.Range("F6").Value = 389
$ws.Range("F10").Value = 552
$ws.Range("F12").Value = 9779
$ws.Range("F14").Value = 51
$ws.Range("F16").Value = 2643
$ws.Range("F17").Value = 206
$ws.Range("F18").Value = 2388
$ws.Range("F19").Value = 2655
$ws.Range("F21").Value = 2076
$ws.Range("F23").Value = 79
$ws.Range("F24").Value = 366
$ws.Range("F26").Value = 297
$ws.Range("F27").Value = 61
$ws.Range("F28").Value = 147
$ws.Range("F30").Value = 1286
$ws.Range("F32").Value = 95
$ws.Range("F33").Value = 122
$ws.Range("F35").Value = 1676
$ws.Range("F37").Value = 2807
$ws.Range("F38").Value = 989
$ws.Range("F41").Value = 358
$ws.Range("F44").Value = 34
$ws.Range("F45").Value = 1275
$ws.Range("F46").Value = 50
$ws.Range("F49").Value = 155
$ws.Range("F50").Value = 155
